$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect it so we can update the data.
$wasProtected = $ws.ProtectContents
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure note (shared string
# referenced by cell A80) from 2021-05-25 to 2021-05-26.
$disclosureText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-26 for illustrative purposes only and are subject to change."
$ws.Range("A80").Value = $disclosureText

# Update the Weight (column D) and Percent Change (column E) figures for each
# holding row (rows 2-77) to the refreshed values.
$ws.Cells.Item(2, 4).Value = 0.0749825898102265
$ws.Cells.Item(2, 5).Value = -0.0003940110323090718
$ws.Cells.Item(3, 4).Value = 0.04585012745713136
$ws.Cells.Item(3, 5).Value = 0.001874779460272036
$ws.Cells.Item(4, 4).Value = 0.03895470610104176
$ws.Cells.Item(4, 5).Value = -0.000913713650087411
$ws.Cells.Item(5, 4).Value = 0.03544178766590093
$ws.Cells.Item(5, 5).Value = -0.008514540214828559
$ws.Cells.Item(6, 4).Value = 0.03324216893408568
$ws.Cells.Item(6, 5).Value = 0.00738085463863869
$ws.Cells.Item(7, 4).Value = 0.03035994386477331
$ws.Cells.Item(7, 5).Value = -0.0001235712079084017
$ws.Cells.Item(8, 4).Value = 0.0309770845216476
$ws.Cells.Item(8, 5).Value = 0.0004526462395542108
$ws.Cells.Item(9, 4).Value = 0.02951095340489092
$ws.Cells.Item(9, 5).Value = -0.005938381937911519
$ws.Cells.Item(10, 4).Value = 0.02726514291418349
$ws.Cells.Item(10, 5).Value = -0.002239398974210882
$ws.Cells.Item(11, 4).Value = 0.02654736518666693
$ws.Cells.Item(11, 5).Value = 0.002522522522522497
$ws.Cells.Item(12, 4).Value = 0.02313227011096682
$ws.Cells.Item(12, 5).Value = 0.001702900607367885
$ws.Cells.Item(13, 4).Value = 0.02403481487405212
$ws.Cells.Item(13, 5).Value = 0.0007141156867414011
$ws.Cells.Item(14, 4).Value = 0.02074985394142391
$ws.Cells.Item(14, 5).Value = -0.007503282686175394
$ws.Cells.Item(15, 4).Value = 0.01999592002496738
$ws.Cells.Item(15, 5).Value = 0.008442855936734883
$ws.Cells.Item(16, 4).Value = 0.01838291318348955
$ws.Cells.Item(16, 5).Value = 0.001020408163265207
$ws.Cells.Item(17, 4).Value = 0.01735515827193507
$ws.Cells.Item(17, 5).Value = -0.001194323450892321
$ws.Cells.Item(18, 4).Value = 0.01738639046627233
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(19, 4).Value = 0.01690895057055522
$ws.Cells.Item(19, 5).Value = -0.0003965953811890799
$ws.Cells.Item(20, 4).Value = 0.01529983602957288
$ws.Cells.Item(20, 5).Value = 0.01167181599725375
$ws.Cells.Item(21, 4).Value = 0.01370502452053177
$ws.Cells.Item(21, 5).Value = 0.001016260162601812
$ws.Cells.Item(22, 4).Value = 0.01557970650375417
$ws.Cells.Item(22, 5).Value = 0.001222067039106101
$ws.Cells.Item(23, 4).Value = 0.01377668036579803
$ws.Cells.Item(23, 5).Value = -0.01244485594466549
$ws.Cells.Item(24, 4).Value = 0.01549604549069461
$ws.Cells.Item(24, 5).Value = 0.0003873623048056274
$ws.Cells.Item(25, 4).Value = 0.01336864534638593
$ws.Cells.Item(25, 5).Value = 0.001831097079715738
$ws.Cells.Item(26, 4).Value = 0.01066602884210612
$ws.Cells.Item(26, 5).Value = 0.0124338298658131
$ws.Cells.Item(27, 4).Value = 0.0120121082809298
$ws.Cells.Item(27, 5).Value = -0.003638521479769619
$ws.Cells.Item(28, 4).Value = 0.01191597314821002
$ws.Cells.Item(28, 5).Value = -0.0004565167769915446
$ws.Cells.Item(29, 4).Value = 0.01193355884321974
$ws.Cells.Item(29, 5).Value = 0.006012425679738254
$ws.Cells.Item(30, 4).Value = 0.01174161684408703
$ws.Cells.Item(30, 5).Value = 0.008547008547008517
$ws.Cells.Item(31, 4).Value = 0.01040102414210637
$ws.Cells.Item(31, 5).Value = 0.0008791981712679142
$ws.Cells.Item(32, 4).Value = 0.01201548473437167
$ws.Cells.Item(32, 5).Value = 0.02064631956912022
$ws.Cells.Item(33, 4).Value = 0.0108594245920263
$ws.Cells.Item(33, 5).Value = -0.002124645892351285
$ws.Cells.Item(34, 4).Value = 0.01079142657132206
$ws.Cells.Item(34, 5).Value = 0.004380361379813902
$ws.Cells.Item(35, 4).Value = 0.01090641356909226
$ws.Cells.Item(35, 5).Value = -0.01096444081351833
$ws.Cells.Item(36, 4).Value = 0.009994302187921668
$ws.Cells.Item(36, 5).Value = -0.01393581081081097
$ws.Cells.Item(37, 4).Value = 0.01004565241735004
$ws.Cells.Item(37, 5).Value = 0.001190392829633735
$ws.Cells.Item(38, 4).Value = 0.008507115132340025
$ws.Cells.Item(38, 5).Value = 0.02388000463047169
$ws.Cells.Item(39, 4).Value = 0.01032857107866634
$ws.Cells.Item(39, 5).Value = -0.01042461225527591
$ws.Cells.Item(40, 4).Value = 0.009254905779339831
$ws.Cells.Item(40, 5).Value = 0.002406854722248886
$ws.Cells.Item(41, 4).Value = 0.008620085637082438
$ws.Cells.Item(41, 5).Value = 0.002067284675980385
$ws.Cells.Item(42, 4).Value = 0.008781123708184742
$ws.Cells.Item(42, 5).Value = -0.004934579439252351
$ws.Cells.Item(43, 4).Value = 0.009611543674136866
$ws.Cells.Item(43, 5).Value = 0.005454776100469516
$ws.Cells.Item(44, 4).Value = 0.008818546067165417
$ws.Cells.Item(44, 5).Value = 0.01913341274568192
$ws.Cells.Item(45, 4).Value = 0.008704121811635534
$ws.Cells.Item(45, 5).Value = -0.009331494332140866
$ws.Cells.Item(46, 4).Value = 0.009287497934091171
$ws.Cells.Item(46, 5).Value = 0.009936984973339724
$ws.Cells.Item(47, 4).Value = 0.008722504724819023
$ws.Cells.Item(47, 5).Value = -0.003483870967741942
$ws.Cells.Item(48, 4).Value = 0.008393581885357299
$ws.Cells.Item(48, 5).Value = 0.008224106913389928
$ws.Cells.Item(49, 4).Value = 0.007922707315777149
$ws.Cells.Item(49, 5).Value = 0.001953298410725557
$ws.Cells.Item(50, 4).Value = 0.009151126731189162
$ws.Cells.Item(50, 5).Value = -0.01555806087936873
$ws.Cells.Item(51, 4).Value = 0.007785491999514673
$ws.Cells.Item(51, 5).Value = 0.002276847830670059
$ws.Cells.Item(52, 4).Value = 0.008199810973943577
$ws.Cells.Item(52, 5).Value = 0.002813776064602536
$ws.Cells.Item(53, 4).Value = 0.006565795088827481
$ws.Cells.Item(53, 5).Value = 0.01757017355903145
$ws.Cells.Item(54, 4).Value = 0.007473967274315895
$ws.Cells.Item(54, 5).Value = 0.008884650135528638
$ws.Cells.Item(55, 4).Value = 0.006555947099622039
$ws.Cells.Item(55, 5).Value = 0.01169527896995715
$ws.Cells.Item(56, 4).Value = 0.006780298362277989
$ws.Cells.Item(56, 5).Value = 0.003948566958607458
$ws.Cells.Item(57, 4).Value = 0.007934290426890216
$ws.Cells.Item(57, 5).Value = 0.01096978580547536
$ws.Cells.Item(58, 4).Value = 0.006552570646180174
$ws.Cells.Item(58, 5).Value = 0.001374098248024591
$ws.Cells.Item(59, 4).Value = 0.006523308049684007
$ws.Cells.Item(59, 5).Value = 0.003105590062111752
$ws.Cells.Item(60, 4).Value = 0.005838263162478866
$ws.Cells.Item(60, 5).Value = 0.002441845521141239
$ws.Cells.Item(61, 4).Value = 0.005668643272211821
$ws.Cells.Item(61, 5).Value = 0.00955500955500943
$ws.Cells.Item(62, 4).Value = 0.005804686208806982
$ws.Cells.Item(62, 5).Value = -0.0009694619486185774
$ws.Cells.Item(63, 4).Value = 0.004955883329671367
$ws.Cells.Item(63, 5).Value = -0.008099924299772798
$ws.Cells.Item(64, 4).Value = 0.005009156261754133
$ws.Cells.Item(64, 5).Value = -0.00943678849610563
$ws.Cells.Item(65, 4).Value = 0.004548645528433043
$ws.Cells.Item(65, 5).Value = -0.00193822425667034
$ws.Cells.Item(66, 4).Value = 0.004513755509533766
$ws.Cells.Item(66, 5).Value = -0.007605036778456586
$ws.Cells.Item(67, 4).Value = 0.004533451487944648
$ws.Cells.Item(67, 5).Value = -0.002358490566037652
$ws.Cells.Item(68, 4).Value = 0.004341415698438552
$ws.Cells.Item(68, 5).Value = 0.006999578729057987
$ws.Cells.Item(69, 4).Value = 0.004144409019143043
$ws.Cells.Item(69, 5).Value = -0.006879695845025791
$ws.Cells.Item(70, 4).Value = 0.003464428812100699
$ws.Cells.Item(70, 5).Value = 0.0160810005955927
$ws.Cells.Item(71, 4).Value = 0.003638832011410387
$ws.Cells.Item(71, 5).Value = 0.01051614150396296
$ws.Cells.Item(72, 4).Value = 0.002935216630274986
$ws.Cells.Item(72, 5).Value = 0.003339138214759396
$ws.Cells.Item(73, 4).Value = 0.002351043289645575
$ws.Cells.Item(73, 5).Value = 0.002034547412933518
$ws.Cells.Item(74, 4).Value = 0.002368582089468598
$ws.Cells.Item(74, 5).Value = 0.003761780312029783
$ws.Cells.Item(75, 4).Value = 0.001876370209943327
$ws.Cells.Item(75, 5).Value = 0.02764170748775374
$ws.Cells.Item(76, 4).Value = 0.001945212344008027
$ws.Cells.Item(76, 5).Value = 0.01610414657666337
$ws.Cells.Item(77, 5).Value = 0.001029626029503738

# Restore sheet protection to match the original workbook state.
if ($wasProtected) {
  $ws.Protect("D382", $true, $true, $true)
}
